# Update the timesheet template:
# - Shift the week-of dates on row 5 (B5:H5) forward by 11 weeks (78 days)
# - Move the active selection from A9 to A12

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B5").Value = 46048
$ws.Range("C5").Value = 46049
$ws.Range("D5").Value = 46050
$ws.Range("E5").Value = 46051
$ws.Range("F5").Value = 46052
$ws.Range("G5").Value = 46053
$ws.Range("H5").Value = 46054

$ws.Range("A12").Select()
